# Fruta / hortaliza, semanal
# The weekly refresh reshuffles the per-row "Fecha" + quality/price/origin
# columns (D, L, M, N, O, P, Q, R, S, T) across the existing data rows
# (rows 2-25). Columns A, B, C, E, F, G, H, I, J, K stay put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# destination row -> source row (values to copy FROM source INTO destination)
$mapping = @{
    2  = 24
    3  = 15
    4  = 14
    5  = 20
    6  = 3
    7  = 2
    8  = 16
    9  = 23
    10 = 9
    11 = 18
    12 = 4
    13 = 5
    14 = 10
    15 = 6
    16 = 21
    17 = 17
    18 = 19
    19 = 8
    20 = 13
    21 = 22
    22 = 7
    23 = 11
    24 = 12
    25 = 25
}

$cols = @(4, 12, 13, 14, 15, 16, 17, 18, 19, 20)   # D, L, M, N, O, P, Q, R, S, T

# Snapshot all the "before" values first so the row-to-row copy uses the
# original data regardless of write order.
$snapshot = @{}
foreach ($row in 2..25) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Cells.Item($row, $col).Value2
    }
    $snapshot[$row] = $rowVals
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Cells.Item($destRow, $col).Value = $srcVals[$col]
    }
}
